$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the label in C2 from "همکف و ویلایی" to "همکف"
$ws.Range("C2").Value = "همکف"

# Reflect the last-used selection cell recorded in the saved sheet view
$ws.Range("J7").Select()
